# Add season-record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties"
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style used by the rest of row 1
# (copy the format from the neighboring "Unnamed: 28" header cell).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-47: constant season record for every player.
$lastRow = 47
$ws.Range("AD2:AD$lastRow").Value = 77
$ws.Range("AE2:AE$lastRow").Value = 85
$ws.Range("AF2:AF$lastRow").Value = 0
